# Update Release-Notes.xlsx - Folder inventory updated on Wed Jun 11 17:17:21 UTC 2025
#
# A new folder run placed "Hackathon - Activate GenAI with Azure" back on top
# of the "Folder Inventory" sheet (its Last Updated / File Count changed), so
# every other row that used to sit above it in the list shifts down by one
# row. The Metadata and Summary sheets are refreshed to match.

$wb = $excel.ActiveWorkbook

$inventory = $wb.Worksheets.Item("Folder Inventory")
$metadata  = $wb.Worksheets.Item("Metadata")
$summary   = $wb.Worksheets.Item("Summary")

# Rows 2-15 (Folder Path, Folder Name, Last Updated, File Count), in the new
# top-to-bottom order. "Hackathon - Activate GenAI with Azure" now leads with
# a refreshed timestamp and file count; the rest retain their previous values
# but are pushed down by one row.
$rows = @(
    @("Hackathon - Activate GenAI with Azure", "2025-06-11 22:47:04 +0530", 2),
    @("Azure Landing Zone", "2025-06-11 20:16:49 +0530", 1),
    @("Microsoft Azure AI Agents", "2025-06-11 20:13:48 +0530", 1),
    @("Azure Local Hands-on Lab", "2025-06-11 19:56:28 +0530", 1),
    @("Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals", "2025-06-11 16:58:40 +0530", 1),
    @("Enterprise-Class Networking in Azure", "2025-06-11 15:50:29 +0530", 1),
    @("Get data into Fabric Lakehouse", "2025-06-11 15:00:50 +0000", 1),
    @("Low Code for Pro-Dev in a Day", "2025-06-11 00:35:20 +0530", 1),
    @("Developing AI Applications with Azure AI Foundry", "2025-06-11 00:33:06 +0530", 1),
    @("Cloud-Native Applications", "2025-06-11 00:33:02 +0530", 1),
    @("Develop Generative AI solutions with Azure OpenAI Service", "2025-06-10 23:22:30 +0530", 1),
    @("Advanced Workflow Automation with GitHub Actions ", "2025-06-10 23:10:36 +0530", 1),
    @("Get Started With OpenAI And Build Natural Language Solution", "2025-06-10 22:51:47 +0530", 1),
    @("Lunch and Learn: Building and Evaluating Prompt Flows with Azure AI Foundry", "2025-06-10 22:48:16 +0530", 1)
)

$r = 2
foreach ($item in $rows) {
    $name = $item[0]
    $lastUpdated = $item[1]
    $fileCount = $item[2]

    $inventory.Cells.Item($r, 1).Value = $name
    $inventory.Cells.Item($r, 2).Value = $name
    $inventory.Cells.Item($r, 3).Value = $lastUpdated
    $inventory.Cells.Item($r, 4).Value = $fileCount
    $inventory.Cells.Item($r, 5).Value = "Root"

    $r = $r + 1
}

# Metadata sheet: refreshed generation timestamp and workflow run number.
# "Workflow Run" is stored as text (e.g. "3"), not a number, so force text
# using a quote-prefixed value and then restore the default (unstyled) cell
# style so no visible number formatting is left behind.
$metadata.Range("B3").Value = "2025-06-11 17:17:20 UTC"
$metadata.Range("B5").Value = "'4"
$metadata.Range("B5").Style = "Normal"

# Summary sheet: most recent update now reflects the new top folder.
$summary.Range("B5").Value = "2025-06-11 22:47:04 +0530"
